$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131, shifting existing rows 131..220 down to 132..221.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new data.
$ws.Cells.Item(131, 1).Value = 11
$ws.Cells.Item(131, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(131, 3).Value = "Bíobío"
$ws.Cells.Item(131, 4).Value = 45126
$ws.Cells.Item(131, 5).Value = 8
$ws.Cells.Item(131, 6).Value = 100112021
$ws.Cells.Item(131, 7).Value = "Ají"
$ws.Cells.Item(131, 8).Value = "Inferno"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 10
$ws.Cells.Item(131, 11).Value = 13000
$ws.Cells.Item(131, 12).Value = 13000
$ws.Cells.Item(131, 13).Value = 13000
$ws.Cells.Item(131, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(131, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(131, 16).Value = 1300
$ws.Cells.Item(131, 17).Value = 10
$ws.Cells.Item(131, 18).Value = "Hortaliza"
